# Bill template update:
#  - {total} placeholder (first occurrence, the "Total" row) -> {finalAmount}
#    (with trailing spacing runs preserved/adjusted)
#  - "Tax rate" label -> "Total GST"
#  - {taxrate} placeholder -> {totalGST}
#  - lone "0" value (under the "Tax" row) -> {totalGST}

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) {total} -> {finalAmount}   (the first "{total}" cell, to the right of
#    {price}; it originally held an 8-space run followed by the {total} run)
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("{total}")
if ($found) {
    $xml = $pkgOpen + '<w:p w14:paraId="0A7A63FF" w14:textId="4942552B" w:rsidR="00885E8B" w:rsidRPr="0093572E" w:rsidRDefault="0079660D" w:rsidP="0079660D"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing w:val="0"/><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>finalAmount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">     </w:t></w:r></w:p>' + $pkgClose
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) "Tax rate" label -> "Total GST"  (kept bold/black formatting, split into
#    "T" + "otal GST")
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("Tax rate")
if ($found) {
    $xml = $pkgOpen + '<w:p w14:paraId="594F66AA" w14:textId="610AACB3" w:rsidR="00885E8B" w:rsidRPr="0093572E" w:rsidRDefault="00885E8B" w:rsidP="003B7694"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing w:val="0"/><w:jc w:val="right"/><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="0093572E"><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>T</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>otal GST</w:t></w:r></w:p>' + $pkgClose
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) {taxrate} -> {totalGST}  (split "taxrate" into "t" + "otalGST", keeping
#    the existing spellStart/spellEnd proofErr markers, and the surrounding
#    "{" / "}" runs, around them)
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("taxrate")
if ($found) {
    $xml = $pkgOpen + '<w:p w14:paraId="00C79398" w14:textId="0E782AD2" w:rsidR="00885E8B" w:rsidRPr="0093572E" w:rsidRDefault="00C075E3" w:rsidP="001E7B62"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing w:val="0"/><w:jc w:val="center"/><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>otalGST</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r></w:p>' + $pkgClose
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 4) lone "0" value (below the "Tax" row) -> {totalGST}
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("0", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $xml = $pkgOpen + '<w:p w14:paraId="286BA3C7" w14:textId="27797B5E" w:rsidR="00885E8B" w:rsidRPr="0093572E" w:rsidRDefault="0035603D" w:rsidP="001E7B62"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:contextualSpacing w:val="0"/><w:jc w:val="center"/><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>totalGST</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="414042"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r></w:p>' + $pkgClose
    $range.InsertXML($xml)
}
